$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.946.74"
$ws.Range("E2").Value = "  -0.56%  "

# Row 3
$ws.Range("D3").Value = "1.638.70"
$ws.Range("E3").Value = "  -0.04%  "

# Row 4
$ws.Range("E4").Value = "  +0.59%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.89"
$ws.Range("E5").Value = "  -0.44%  "

# Row 7
$ws.Range("E7").Value = "  +0.61%  "

# Row 8
$ws.Range("E8").Value = "  -0.72%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0639"
$ws.Range("E9").Value = "  +0.64%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.65"
$ws.Range("E10").Value = "  -1.01%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0795"

# Row 12
$ws.Range("D12").Value = "1.866.04"
$ws.Range("E12").Value = "  +0.06%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.26"
$ws.Range("E13").Value = "  -0.27%  "

# Row 14
$ws.Range("D14").Value = "1.654.83"
$ws.Range("E14").Value = "  +1.26%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.545"
$ws.Range("E15").Value = "  -1.47%  "

# Row 16
$ws.Range("E16").Value = "  -0.43%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.69"
$ws.Range("E17").Value = "  -1.04%  "

# Row 18
$ws.Range("D18").Value = "25.957.77"
$ws.Range("E18").Value = "  -0.40%  "

# Row 19
$ws.Range("E19").Value = "  +0.61%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.07"
$ws.Range("E20").Value = "  +0.31%  "

# Row 21
$ws.Range("E21").Value = "  -1.75%  "

# Row 22
$ws.Range("E22").Value = "  -0.96%  "

# Row 23
$ws.Range("E23").Value = "  -1.16%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.14"
$ws.Range("E24").Value = "  +1.06%  "

# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.78"
$ws.Range("E25").Value = "  -0.44%  "

# Row 26
$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.69%  "

# Row 27
$ws.Range("E27").Value = "  +1.71%  "

# Row 28
$ws.Range("E28").Value = "  -1.00%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.51"
$ws.Range("E29").Value = "  -0.48%  "

# Row 30
$ws.Range("E30").Value = "  -0.22%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0502"
$ws.Range("E31").Value = "  +0.48%  "

# Row 32
$ws.Range("E32").Value = "  -1.55%  "

# Row 33
$ws.Range("E33").Value = "  -0.43%  "

# Row 34
$ws.Range("E34").Value = "  -3.15%  "

# Row 35
$ws.Range("E35").Value = "  +1.44%  "

# Row 36
$ws.Range("E36").Value = "  -0.47%  "

# Row 37
$ws.Range("D37").Value = "1.138.63"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.544"
$ws.Range("E38").Value = "  -1.07%  "

# Row 39
$ws.Range("E39").Value = "  -2.00%  "

# Row 40
$ws.Range("E40").Value = "  +0.10%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.36"
$ws.Range("E41").Value = "  -1.30%  "

# Row 42
$ws.Range("E42").Value = "  +0.85%  "

# Row 43
$ws.Range("E43").Value = "  -3.00%  "

# Row 44
$ws.Range("D44").Value = "1.776.08"
$ws.Range("E44").Value = "  +0.12%  "

# Row 45
$ws.Range("E45").Value = "  +11.38%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.56"
$ws.Range("E46").Value = "  +1.07%  "

# Row 47
$ws.Range("E47").Value = "  +2.60%  "

# Row 48
$ws.Range("E48").Value = "  -1.20%  "

# Row 49
$ws.Range("E49").Value = "  -0.73%  "

# Row 50
$ws.Range("E50").Value = "  -0.46%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0965"
$ws.Range("E51").Value = "  -1.06%  "
